$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 values (B:G)
$ws.Range("B2").Value = 0.08780261824621656
$ws.Range("C2").Value = 0.7725920269955064
$ws.Range("D2").Value = 2.162285864834432
$ws.Range("E2").Value = 1.470471307042212
$ws.Range("F2").Value = 1.484067398835772
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.2407776619532973
$ws.Range("C3").Value = 0.7527870312651909
$ws.Range("D3").Value = 2.110441622564052
$ws.Range("E3").Value = 1.452735909435728
$ws.Range("F3").Value = 1.44883219668836
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.1187069702962777
$ws.Range("C4").Value = 0.783804746145044
$ws.Range("D4").Value = 2.289524756223591
$ws.Range("E4").Value = 1.513117561930861
$ws.Range("F4").Value = 1.525893334249272
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.2312605598051813
$ws.Range("C5").Value = 0.8154870963700255
$ws.Range("D5").Value = 2.323257756554268
$ws.Range("E5").Value = 1.524223656998627
$ws.Range("F5").Value = 1.524407623678418
$ws.Range("G5").Value = 43

# Add new rows 6-11 with formatting copied from row 5 (A col) and data cells
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = "Q4"
$ws.Range("B6").Value = 0.1543123306983461
$ws.Range("C6").Value = 0.8126759977456647
$ws.Range("D6").Value = 2.40121625977409
$ws.Range("E6").Value = 1.549585834916572
$ws.Range("F6").Value = 1.5605734288761
$ws.Range("G6").Value = 42

$ws.Range("A5").Copy($ws.Range("A7"))
$ws.Range("A7").Value = "Q5"
$ws.Range("B7").Value = 0.2357582618036748
$ws.Range("C7").Value = 0.8352230691908547
$ws.Range("D7").Value = 2.420446362568442
$ws.Range("E7").Value = 1.555778378358705
$ws.Range("F7").Value = 1.556915545132023
$ws.Range("G7").Value = 41

$ws.Range("A5").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "Q6"
$ws.Range("B8").Value = 0.1764347851107931
$ws.Range("C8").Value = 0.8312978488293193
$ws.Range("D8").Value = 2.481252276672181
$ws.Range("E8").Value = 1.575199122864211
$ws.Range("F8").Value = 1.585227652752555
$ws.Range("G8").Value = 40

$ws.Range("A5").Copy($ws.Range("A9"))
$ws.Range("A9").Value = "Q7"
$ws.Range("B9").Value = 0.2294621647362194
$ws.Range("C9").Value = 0.8569143726678996
$ws.Range("D9").Value = 2.527131915378803
$ws.Range("E9").Value = 1.589695541724516
$ws.Range("F9").Value = 1.593611276175175
$ws.Range("G9").Value = 39

$ws.Range("A5").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.1874282125274516
$ws.Range("C10").Value = 0.8425231446840433
$ws.Range("D10").Value = 2.595659384232942
$ws.Range("E10").Value = 1.611105019616332
$ws.Range("F10").Value = 1.621645326275095
$ws.Range("G10").Value = 38

$ws.Range("A5").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "Q9"
$ws.Range("B11").Value = 0.1764013765962084
$ws.Range("C11").Value = 0.8338957443602626
$ws.Range("D11").Value = 2.565872115365881
$ws.Range("E11").Value = 1.60183398495783
$ws.Range("F11").Value = 1.614052205362943
$ws.Range("G11").Value = 37
